$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (fra/French) - native_name column gets the mis-encoded value, is_active -> boolean
$ws.Range("D3").Value = "franÃ§ais"

# Row 4 (ara/Arabic) - family and native_name columns get mis-encoded values
$ws.Range("C4").Value = "Ø§Ù„Ù‡Ù†Ø¯Ùˆ Ø£ÙˆØ±ÙˆØ¨ÙŠØ©"
$ws.Range("D4").Value = "Ø¹Ø±Ø¨ÙŠ"

# is_active column (E) becomes a native boolean TRUE instead of the text "TRUE" shared string
$ws.Range("E2:E4").Value = $true

# Remove the direct cell formatting/styles from A2:E4 (rows lose their s="4"/s="1" attributes)
$ws.Range("A2:E4").Style = "Normal"

# Remove the kan/hin/tam rows (rows 5-7), leaving only eng/fra/ara
$ws.Rows("5:7").Delete() | Out-Null

# Restore the selection to match the saved view state
$ws.Range("B8").Select() | Out-Null
